# "edit report paging and searching"
#
# - B1 gets a date/time value (44542.375) that was previously blank
#   (cell already carries the date-time number format, style index 1).
# - The sheet view zoom is bumped to 145% and the active selection moves
#   to D3 (was C2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate B1 with the missing date/time value.
$ws.Range("B1").Value = 44542.375

# Zoom the sheet view to 145% and move the selection to D3.
$excel.ActiveWindow.Zoom = 145
$ws.Range("D3").Select()
